$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before the old "comments" column (J), shifting
# J:N -> L:P. This updates dimension/col widths/row spans automatically.
$ws.Columns("J:K").Insert()

# New header cells for the inserted columns.
$ws.Range("J1").Value = "auxillaryHash"
$ws.Range("K1").Value = "auxillaryHash.cell_type"

# New data cells on the "linked_table" settings row (row 7).
$j7 = $ws.Range("J7")
$j7.WrapText = $true
$j7.Value = "'household_id='+escape(data('household_id'))"
$ws.Range("K7").Value = "formula"

# Update the comments cell (now L7) with the combined auxillaryHash +
# joined_through_name explanation (replacing the old joined_through_name-only text).
$l7 = $ws.Range("L7")
$l7.Value = "auxillaryHash defines the auxillary hash to supply when creating a new sub-form. This is an ampersand-separated list of elementName=value pairs that will be used to initialize the subform. The joined_through_name value identifies the name (elementName) in the model that should be used when scanning in the joins lists for the table_id to discover the foreign key column to filter on in the subform.
If this is omitted, we would probably just scan the entire model to see if table_id appears anywhere and use the first match we find. "

$b1 = $l7.Characters(1, 13)
$b1.Font.Bold = $true

$b2 = $l7.Characters(194, 19)
$b2.Font.Bold = $true

# Taller row to fit the longer comment.
$ws.Rows("7:7").RowHeight = 220.5

# Update the view's active selection.
$ws.Range("J8").Select()
